# Generate Report for Handoff
#
# A new handoff cycle was kicked off for the e2e markdown file: it was
# regenerated under a new GUID and its handoff .xlf files now carry a new
# content hash. This refreshes the localization-status workbook to reflect
# that: the File Name / Path / Source File Name columns (and their
# hyperlinks) move to the new GUID, the Latest Handoff File/Datetime columns
# pick up the new hash + timestamps, and the stale Latest Target
# File / Latest Handback File / Latest Handback DateTime columns (left over
# from the previous handback) are cleared out for the new cycle.

$wb = $excel.ActiveWorkbook

$oldGuid = "9db4f506-09c0-41b9-b5f5-b7ddff47b76a"
$newGuid = "b7edaa1b-dcd9-4334-bd51-d638849d9e7e"

$oldHash = "f97d646ee18c8449b3c72bfcb12479a71bb09980"
$newHash = "d720e9e72f739d8a57b61a9ed74066cc946926fc"

$newOverviewDate  = "2016-08-28 11:08:17"
$newZhHandoffDate = "2016-08-28 11:08:12"
$newDeHandoffDate = "2016-08-28 11:08:17"

$clearedHandbackDate = "0001-01-01 00:00:00"

$hyperlinkFontColor = 15570276   # BGR for RGB FF6495ED, matching the workbook's custom HyperLink style

function Restore-HyperlinkLook($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkFontColor
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newOverviewDate

$overviewHyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/$oldGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewHyperlinkAddr, "", "", "e2e\$newGuid.md") | Out-Null
Restore-HyperlinkLook $wsOverview.Range("B2")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $clearedHandbackDate
$wsZh.Range("I2").Style = "Normal"

$zhHyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/$oldGuid.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhHyperlinkAddr, "", "", "$newGuid.md") | Out-Null
Restore-HyperlinkLook $wsZh.Range("A2")

$wsZh.Columns.Item(9).AutoFit() | Out-Null
$wsZh.Columns.Item(10).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newDeHandoffDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $clearedHandbackDate
$wsDe.Range("I2").Style = "Normal"

$deHyperlinkAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a3f364f6003c6ee15ff7d8e3199207c0783deec7/e2e/$oldGuid.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deHyperlinkAddr, "", "", "$newGuid.md") | Out-Null
Restore-HyperlinkLook $wsDe.Range("A2")

$wsDe.Columns.Item(9).AutoFit() | Out-Null
$wsDe.Columns.Item(10).AutoFit() | Out-Null
